# Refresh the crypto "Price" column (D) with newly scraped values.
#
# The source values are stored as text (inline strings), not numbers, so a
# plain assignment like `$ws.Range("D2").Value = "270.55"` would get
# re-typed by Excel into a numeric cell. To keep the values as text we
# prefix the literal with a leading apostrophe (the same trick you'd use
# typing into the grid), then reset the cell style back to "Normal" so the
# quote-prefix formatting picked up along the way doesn't linger on the
# cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $newValue
    $rng.Style = "Normal"
}

Set-TextValue "D2"  "270.55"
Set-TextValue "D4"  "6.327"
Set-TextValue "D5"  "0.06285"
Set-TextValue "D6"  "3.566"
Set-TextValue "D7"  "6.574"
Set-TextValue "D8"  "1.378"
Set-TextValue "D9"  "0.8279"
Set-TextValue "D10" "0.01383"
Set-TextValue "D11" "0.1579"
Set-TextValue "D12" "0.08333"
Set-TextValue "D13" "0.03415"
Set-TextValue "D14" "0.03210"
Set-TextValue "D15" "4.074"
Set-TextValue "D16" "0.09266"
Set-TextValue "D17" "0.001658"
Set-TextValue "D18" "0.04692"
Set-TextValue "D19" "0.006313"
Set-TextValue "D21" "0.001064"
Set-TextValue "D22" "0.0001493"
Set-TextValue "D23" "3.735"
Set-TextValue "D24" "2.326"
Set-TextValue "D25" "0.3332"
Set-TextValue "D26" "0.1255"
Set-TextValue "D28" "0.0002713"
Set-TextValue "D40" "0.04721"
Set-TextValue "D41" "0.007062"
Set-TextValue "D42" "0.1171"
Set-TextValue "D43" "0.003638"
Set-TextValue "D44" "0.01171"
Set-TextValue "D45" "0.00005992"
Set-TextValue "D46" "0.0009841"
Set-TextValue "D47" "0.00000000746"
Set-TextValue "D48" "0.7782"
Set-TextValue "D49" "0.002404"
Set-TextValue "D50" "0.00001294"
Set-TextValue "D51" "0.01234"
